# The source workbook re-ran its word-frequency tally over an updated corpus.
# Every word (column B) and its count (column C) for rows 2-101 changes, and the
# word list itself is re-ranked, so every data row is rewritten from scratch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data body first (keeps the header row/styles untouched) so the
# saved shared-string table reflects only the strings actually used below.
$ws.Range("A2:C101").ClearContents()

$words = @(
    "식당", "없", "메뉴", "밥", "아침", "배식", "식사", "좋", "사람", "중앙", "빵", "이용", "복지관", "시간", "음식", "식", "식단",
    "분", "오늘", "많", "맛있", "맛", "글", "생각", "편", "반찬", "의식", "개선", "감사", "우유", "층", "식권", "부탁", "주세",
    "관련", "답변", "말", "점심", "청운", "라면", "직원", "제공", "안", "조식", "중", "경우", "그렇", "국", "안녕", "업체", "고객",
    "전", "중식", "정도", "청주", "문의", "이렇", "건의", "데", "아주머니", "바나나", "앞", "신경", "확인", "기분", "줄", "기호", "김치",
    "요즘", "회사", "불편", "일", "하세", "양", "문제", "영양사", "운영", "사항", "어떻", "기숙사", "수고", "등", "번", "후", "코너",
    "하이닉스", "한식", "야식", "석식", "저녁", "고담", "종류", "시", "준비", "추가", "금일", "안되", "주말", "인", "간"
)

$counts = @(
    1917, 871, 835, 626, 577, 564, 546, 539, 523, 522, 508, 464, 450, 446, 435, 430, 428, 409, 408, 389,
    378, 374, 371, 354, 351, 335, 326, 326, 311, 298, 281, 264, 255, 236, 227, 219, 218, 216, 216, 209,
    207, 202, 201, 199, 193, 193, 189, 189, 185, 185, 184, 184, 183, 178, 176, 175, 174, 171, 159, 158,
    156, 156, 155, 155, 154, 151, 150, 149, 147, 147, 146, 141, 139, 137, 136, 135, 133, 131, 129, 127,
    126, 126, 126, 125, 124, 122, 120, 120, 119, 119, 119, 118, 116, 116, 115, 113, 111, 109, 108, 108
)

for ($i = 0; $i -lt $words.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $i
    $ws.Range("B$row").Value = $words[$i]
    $ws.Range("C$row").Value = $counts[$i]
}

